$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are stored as text (matching original inlineStr type)
# so values like "246.31" or "24.01" are not reinterpreted as floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.31'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '2'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '24.01'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '2'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.365'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '2'

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '2'

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.465'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5KuCoinTokenKCS'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '2'

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.341'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '6GateTokenGT'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '2'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8084'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '2'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9193'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '2'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1397'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '2'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07356'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '2'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03176'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '2'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03060'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '2'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09377'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '2'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.863'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '2'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001552'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '2'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04708'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '2'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0006004'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '2'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.005933'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '2'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001276'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '2'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.004663'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '2'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.00008806'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '2'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.620'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '2'

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '2'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3180'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '2'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1319'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '2'

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '2'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002351'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '27UpBotsUBXTWorstin24h'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '2'

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '2'

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '2'

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '2'

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '2'

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '2'

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '2'

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '2'

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '2'

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '2'

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '2'

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '2'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03836'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '2'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.004903'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '2'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006404'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '41KickTokenKICKBestin24h'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '2'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1064'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '2'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007799'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '2'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005312'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '2'

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '2'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6860'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '2'

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '2'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '2'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '2'

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '2'
